$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.200.12'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '1.909.68'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7383'
$ws.Range('E5').Value = '  -4.29%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '243.65'
$ws.Range('E6').Value = '  -2.15%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3129'
$ws.Range('E8').Value = '  -2.67%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '26.96'
$ws.Range('E9').Value = '  -4.08%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06959'
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.7789'
$ws.Range('E11').Value = '  -1.21%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07971'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').Value = '1.917.65'
$ws.Range('E13').Value = '  -0.80%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.281'
$ws.Range('E14').Value = '  -1.96%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '91.66'
$ws.Range('E15').Value = '  -3.51%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '30.262.60'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.27'
$ws.Range('E17').Value = '  -2.48%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.878'
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '244.31'
$ws.Range('E19').Value = '  -5.36%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007843'
$ws.Range('E20').Value = '  -2.40%  '
$ws.Range('D21').Value = '2.187.31'
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.650'
$ws.Range('E24').Value = '  -2.55%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.423'
$ws.Range('E25').Value = '  -2.09%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '165.47'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.98'
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.1270'
$ws.Range('E28').Value = '  -6.55%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.110'
$ws.Range('E29').Value = '  -8.52%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.548'
$ws.Range('E30').Value = '  +1.14%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.347'
$ws.Range('E31').Value = '  -1.79%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.311'
$ws.Range('E32').Value = '  -3.04%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.074'
$ws.Range('E33').Value = '  -2.46%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05178'
$ws.Range('E34').Value = '  -0.78%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.294'
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7455'
$ws.Range('E36').Value = '  -1.03%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.764'
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01940'
$ws.Range('E38').Value = '  -1.91%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.792'
$ws.Range('E39').Value = '  -0.81%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.377'
$ws.Range('E40').Value = '  -1.57%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '75.38'
$ws.Range('E41').Value = '  -4.01%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.4479'
$ws.Range('E42').Value = '  -1.15%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.935'
$ws.Range('E43').Value = '  -2.70%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.8350'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.662'
$ws.Range('E46').Value = '  +1.07%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '101.36'
$ws.Range('E47').Value = '  -0.83%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.843'
$ws.Range('E48').Value = '  -0.46%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '37.50'
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '943.18'
$ws.Range('E51').Value = '  -4.34%  '
